# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) columns hold text-formatted values (e.g.
# "3.659.38", "  +7.45%  ") rather than numbers, so several D-column updates
# use a leading quote-prefix to force text entry, then reset the cell style
# back to "Normal" so no stray number-format style is left attached to the
# cell (matches the original workbook, which has no explicit style on these
# cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.151.44'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '3.659.38'
$ws.Range("E3").Value = '  +7.45%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''594.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").Value = '''181.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = '3.648.63'
$ws.Range("E7").Value = '  +7.33%  '
$ws.Range("D8").Value = '''0.608'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '''0.204'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.70%  '
$ws.Range("D11").Value = '''0.606'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("D12").Value = '''50.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.02%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").Value = '''693.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '4.236.99'
$ws.Range("E15").Value = '  +7.13%  '
$ws.Range("D16").Value = '''9.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.04%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '72.195.98'
$ws.Range("E17").Value = '  +3.85%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.654.44'
$ws.Range("E18").Value = '  +6.97%  '
$ws.Range("E19").Value = '  +2.41%  '
$ws.Range("D20").Value = '''18.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.45%  '
$ws.Range("D21").Value = '''11.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").Value = '''0.937'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.07%  '
$ws.Range("D23").Value = '''5.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.95%  '
$ws.Range("D24").Value = '''18.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.69%  '
$ws.Range("D25").Value = '''103.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").Value = '''4.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("D27").Value = '''2.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.22%  '
$ws.Range("D28").Value = '''10.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.75%  '
$ws.Range("D29").Value = '''35.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.75%  '
$ws.Range("D30").Value = '''9.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.69%  '
$ws.Range("E31").Value = '  +5.73%  '
$ws.Range("D32").Value = '''4.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.38%  '
$ws.Range("D33").Value = '''582.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("D34").Value = '''11.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("D36").Value = '''59.64'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '3.678.24'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("D40").Value = '''36.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").Value = '0.0₃0774'
$ws.Range("E41").Value = '  +7.18%  '
$ws.Range("D42").Value = '''3.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.72%  '
$ws.Range("E43").Value = '  +8.59%  '
$ws.Range("D44").Value = '''2.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("D45").Value = '''0.351'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''3.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.12%  '
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '''2.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.73%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.133'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '''1.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.81%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '''1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '''3.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.81%  '
